# Append 3 new rows (20-22) of CSV-imported data to sheet1, matching the
# "working with csv file type" commit: dimension/ignoredErrors ranges grow
# from A1:G19 to A1:G22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("62107f9c1a601c3edefd10d6", "manjit nayak", "manjit@gmail.com", 87655566777, "kolkata", 5, 0),
    @("6210b85eb1075b536f7fbc70", "niyatee", "gudi@gmail.com", 1234567890, "koraput", 1, 0),
    @("622734541f0216479f759b58", "suraj", "suraj@gmail.com", 1123344444, "pune", 1, 0)
)

$startRow = 20
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}

# The original sheet carries a "numbers stored as text" ignored-error marker
# over A1:G19 (left over from the CSV import). Re-assert/extend it over the
# full, newly-grown range so it keeps covering every data row.
$ws.Range("A1:G22").Errors.Item(3).Ignore = $true
